$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 5.770351896557194
$ws.Range("D2").Value = 5.69151669042765
$ws.Range("E2").Value = 11.13952402176482
$ws.Range("F2").Value = 46.52596488021071
$ws.Range("G2").Value = 3.730380938232589
$ws.Range("I2").Value = 33.19093387126759
$ws.Range("J2").Value = 9.827350157133347
$ws.Range("K2").Value = 20.62939283666754
$ws.Range("M2").Value = 20.54564155860964
$ws.Range("N2").Value = 20.60907617808908
$ws.Range("C3").Value = 5.770248028501598
$ws.Range("D3").Value = 5.688085643398176
$ws.Range("E3").Value = 11.15636483053418
$ws.Range("F3").Value = 46.47776743339227
$ws.Range("G3").Value = 3.734196434935924
$ws.Range("I3").Value = 33.16612704358795
$ws.Range("J3").Value = 9.850810723537979
$ws.Range("K3").Value = 20.32242842836381
$ws.Range("M3").Value = 20.44062985909191
$ws.Range("N3").Value = 20.68595807129309
$ws.Range("C4").Value = 5.770238432072141
$ws.Range("D4").Value = 5.686201248938922
$ws.Range("E4").Value = 11.16802994956298
$ws.Range("F4").Value = 46.46062048992555
$ws.Range("G4").Value = 3.736660107638533
$ws.Range("I4").Value = 33.15970658107354
$ws.Range("J4").Value = 9.866251971988234
$ws.Range("K4").Value = 20.13717117331976
$ws.Range("M4").Value = 20.38045079820971
$ws.Range("N4").Value = 20.73519146530589
$ws.Range("C5").Value = 5.770248372426544
$ws.Range("D5").Value = 5.685489611397615
$ws.Range("E5").Value = 11.17311697285297
$ws.Range("F5").Value = 46.45676006599231
$ws.Range("G5").Value = 3.737694602690587
$ws.Range("I5").Value = 33.15930268124468
$ws.Range("J5").Value = 9.872805340024966
$ws.Range("K5").Value = 20.06258313281378
$ws.Range("M5").Value = 20.35702659666847
$ws.Range("N5").Value = 20.75576609386929
$ws.Range("C6").Value = 5.770250866291192
$ws.Range("D6").Value = 5.685374853894723
$ws.Range("E6").Value = 11.17398181260845
$ws.Range("F6").Value = 46.4563077853845
$ws.Range("G6").Value = 3.737868227068657
$ws.Range("I6").Value = 33.15936911861826
$ws.Range("J6").Value = 9.873909291380221
$ws.Range("K6").Value = 20.05025533110415
$ws.Range("M6").Value = 20.35320392346978
$ws.Range("N6").Value = 20.75921344766788
$ws.Range("C7").Value = 5.770238509767752
$ws.Range("D7").Value = 5.686191423204188
$ws.Range("E7").Value = 11.16809720470526
$ws.Range("F7").Value = 46.46055577111017
$ws.Range("G7").Value = 3.736673935439197
$ws.Range("I7").Value = 33.15969218120019
$ws.Range("J7").Value = 9.866339296055253
$ws.Range("K7").Value = 20.1361614596395
$ws.Range("M7").Value = 20.38013041738201
$ws.Range("N7").Value = 20.73546686813984
$ws.Range("C8").Value = 5.770304972966802
$ws.Range("D8").Value = 5.690287754353875
$ws.Range("E8").Value = 11.14505593355799
$ws.Range("F8").Value = 46.5067622822386
$ws.Range("G8").Value = 3.731671490198333
$ws.Range("I8").Value = 33.18055037199951
$ws.Range("J8").Value = 9.835224461325046
$ws.Range("K8").Value = 20.52294119911278
$ws.Range("M8").Value = 20.50855283595229
$ws.Range("N8").Value = 20.63516529471327
$ws.Range("C9").Value = 5.77085545589089
$ws.Range("D9").Value = 5.700068515036304
$ws.Range("E9").Value = 11.11037213366758
$ws.Range("F9").Value = 46.69614033270636
$ws.Range("G9").Value = 3.722816017702566
$ws.Range("I9").Value = 33.29143222207344
$ws.Range("J9").Value = 9.782416735349878
$ws.Range("K9").Value = 21.30269683786376
$ws.Range("M9").Value = 20.79361990741369
$ws.Range("N9").Value = 20.45448002111445
$ws.Range("C10").Value = 5.771503294950709
$ws.Range("D10").Value = 5.708299573669405
$ws.Range("E10").Value = 11.09127537002868
$ws.Range("F10").Value = 46.8953431511068
$ws.Range("G10").Value = 3.716884232055068
$ws.Range("I10").Value = 33.41554633571101
$ws.Range("J10").Value = 9.748602896483662
$ws.Range("K10").Value = 21.88283838264736
$ws.Range("M10").Value = 21.02204999988664
$ws.Range("N10").Value = 20.33137036538443
$ws.Range("C11").Value = 5.771848236876709
$ws.Range("D11").Value = 5.712266682521845
$ws.Range("E11").Value = 11.08397060162637
$ws.Range("F11").Value = 46.99892206300469
$ws.Range("G11").Value = 3.714308829089251
$ws.Range("I11").Value = 33.48123046553361
$ws.Range("J11").Value = 9.734298197834015
$ws.Range("K11").Value = 22.14714245339655
$ws.Range("M11").Value = 21.12979606981303
$ws.Range("N11").Value = 20.27743252112292
$ws.Range("C12").Value = 5.771985865469485
$ws.Range("D12").Value = 5.713800562154587
$ws.Range("E12").Value = 11.08140291599459
$ws.Range("F12").Value = 47.03999645789258
$ws.Range("G12").Value = 3.713351155421246
$ws.Range("I12").Value = 33.50742325549227
$ws.Range("J12").Value = 9.729036006702172
$ws.Range("K12").Value = 22.24718756814027
$ws.Range("M12").Value = 21.171121062036
$ws.Range("N12").Value = 20.25730278275802
$ws.Range("C13").Value = 5.771955916488689
$ws.Range("D13").Value = 5.713468813934925
$ws.Range("E13").Value = 11.08194709135381
$ws.Range("F13").Value = 47.03106824060622
$ws.Range("G13").Value = 3.713556627669037
$ws.Range("I13").Value = 33.50172360530855
$ws.Range("J13").Value = 9.730162437668046
$ws.Range("K13").Value = 22.22564441664806
$ws.Range("M13").Value = 21.16219812729589
$ws.Range("N13").Value = 20.26162497223083
$ws.Range("C14").Value = 5.771859420476643
$ws.Range("D14").Value = 5.712392244342261
$ws.Range("E14").Value = 11.08375538084984
$ws.Range("F14").Value = 47.00226425582086
$ws.Range("G14").Value = 3.714229689095637
$ws.Range("I14").Value = 33.48335895985907
$ws.Range("J14").Value = 9.733862175374993
$ws.Range("K14").Value = 22.15537459878097
$ws.Range("M14").Value = 21.13318556097817
$ws.Range("N14").Value = 20.27577052550578
$ws.Range("C15").Value = 5.77180121994878
$ws.Range("D15").Value = 5.711736920200294
$ws.Range("E15").Value = 11.08488884726154
$ws.Range("F15").Value = 46.98486166441962
$ws.Range("G15").Value = 3.714644244216166
$ws.Range("I15").Value = 33.47228170803268
$ws.Range("J15").Value = 9.736148509482268
$ws.Range("K15").Value = 22.11232410028343
$ws.Range("M15").Value = 21.11548193117335
$ws.Range("N15").Value = 20.28447350036178
$ws.Range("C16").Value = 5.771481744408993
$ws.Range("D16").Value = 5.708044761214263
$ws.Range("E16").Value = 11.09178054516132
$ws.Range("F16").Value = 46.88883382651453
$ws.Range("G16").Value = 3.717055006180158
$ws.Range("I16").Value = 33.41143879704394
$ws.Range("J16").Value = 9.749559401363046
$ws.Range("K16").Value = 21.86556500118136
$ws.Range("M16").Value = 21.01508335236895
$ws.Range("N16").Value = 20.33493674991443
$ws.Range("C17").Value = 5.77129846751706
$ws.Range("D17").Value = 5.705836503956615
$ws.Range("E17").Value = 11.09636223437478
$ws.Range("F17").Value = 46.8332354013106
$ws.Range("G17").Value = 3.718565353587715
$ws.Range("I17").Value = 33.37647189400935
$ws.Range("J17").Value = 9.758062321422511
$ws.Range("K17").Value = 21.71421685765501
$ws.Range("M17").Value = 20.95445454801776
$ws.Range("N17").Value = 20.36642213387529
$ws.Range("C18").Value = 5.771197792790645
$ws.Range("D18").Value = 5.704587364714141
$ws.Range("E18").Value = 11.09912764395415
$ws.Range("F18").Value = 46.80247748208151
$ws.Range("G18").Value = 3.71944564906234
$ws.Range("I18").Value = 33.35722858391662
$ws.Range("J18").Value = 9.763054398746069
$ws.Range("K18").Value = 21.62720825588867
$ws.Range("M18").Value = 20.91994484630035
$ws.Range("N18").Value = 20.3847261926863
$ws.Range("C19").Value = 5.771164527052366
$ws.Range("D19").Value = 5.704168045482725
$ws.Range("E19").Value = 11.10008632573978
$ws.Range("F19").Value = 46.7922733903
$ws.Range("G19").Value = 3.71974569487046
$ws.Range("I19").Value = 33.35086251074743
$ws.Range("J19").Value = 9.764762057159663
$ws.Range("K19").Value = 21.59775902667766
$ws.Range("M19").Value = 20.9083235055636
$ws.Range("N19").Value = 20.3909570880114
$ws.Range("C20").Value = 5.771317488536083
$ws.Range("D20").Value = 5.706069407091765
$ws.Range("E20").Value = 11.09586103867727
$ws.Range("F20").Value = 46.83902767919545
$ws.Range("G20").Value = 3.71840337647644
$ws.Range("I20").Value = 33.38010431376477
$ws.Range("J20").Value = 9.757146675968343
$ws.Range("K20").Value = 21.73032440966023
$ws.Range("M20").Value = 20.96087127620231
$ws.Range("N20").Value = 20.3630503430661
$ws.Range("C21").Value = 5.771887575251793
$ws.Range("D21").Value = 5.712707604019142
$ws.Range("E21").Value = 11.08321885861893
$ws.Range("F21").Value = 47.01067455254809
$ws.Range("G21").Value = 3.714031518527257
$ws.Range("I21").Value = 33.4887173497442
$ws.Range("J21").Value = 9.732771276417029
$ws.Range("K21").Value = 22.17601642812579
$ws.Range("M21").Value = 21.14169325759932
$ws.Range("N21").Value = 20.27160763086982
$ws.Range("C22").Value = 5.772300901098903
$ws.Range("D22").Value = 5.717230171469791
$ws.Range("E22").Value = 11.0761131407298
$ws.Range("F22").Value = 47.13363891227642
$ws.Range("G22").Value = 3.711276649049643
$ws.Range("I22").Value = 33.56738968304488
$ws.Range("J22").Value = 9.71774207631546
$ws.Range("K22").Value = 22.46701919752777
$ws.Range("M22").Value = 21.26291131742122
$ws.Range("N22").Value = 20.21356531817131
$ws.Range("C23").Value = 5.77207664440922
$ws.Range("D23").Value = 5.714799681528981
$ws.Range("E23").Value = 11.07979987284014
$ws.Range("F23").Value = 47.06702864470228
$ws.Range("G23").Value = 3.712737642435974
$ws.Range("I23").Value = 33.52470007838728
$ws.Range("J23").Value = 9.725681030210914
$ws.Range("K23").Value = 22.31176294257294
$ws.Range("M23").Value = 21.19794595782242
$ws.Range("N23").Value = 20.24438668100015
$ws.Range("C24").Value = 5.771308874497367
$ws.Range("D24").Value = 5.705964047958625
$ws.Range("E24").Value = 11.09608722014776
$ws.Range("F24").Value = 46.83640523192545
$ws.Range("G24").Value = 3.718476569015298
$ws.Range("I24").Value = 33.3784594187516
$ws.Range("J24").Value = 9.757560316459168
$ws.Range("K24").Value = 21.72304216715834
$ws.Range("M24").Value = 20.95796919174363
$ws.Range("N24").Value = 20.36457409822556
$ws.Range("C25").Value = 5.770662885560859
$ws.Range("D25").Value = 5.697237531385483
$ws.Range("E25").Value = 11.11863245111938
$ws.Range("F25").Value = 46.63433899334277
$ws.Range("G25").Value = 3.725110263354056
$ws.Range("I25").Value = 33.25393810619119
$ws.Range("J25").Value = 9.795825958408386
$ws.Range("K25").Value = 21.0900672106898
$ws.Range("M25").Value = 20.71306960625162
$ws.Range("N25").Value = 20.50165885061169
